# Daily attendance processing - 2025-12-27 19:52:00
# Normalize the "Recorded By" (column G) values: when the literal token
# "System" is the first entry in a comma-separated list of recorders,
# move it to the end of the list instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $current = $cell.Value2

    if ($current -ne $null -and $current.GetType().Name -eq "String" -and $current -like "*,*") {
        $parts = $current -split ", "
        if ($parts.Count -gt 1 -and $parts[0] -ceq "System") {
            $rest = $parts[1..($parts.Count - 1)]
            $newValue = ($rest + $parts[0]) -join ", "
            $cell.Value2 = $newValue
        }
    }
}
